$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple price updates (column D). Leading "'" forces text storage,   ---
# --- matching the source data (prices are stored as text, not numbers). ---
$ws.Range("D2").Value  = "'244.73"
$ws.Range("D3").Value  = "'23.04"
$ws.Range("D4").Value  = "'5.414"
$ws.Range("D5").Value  = "'0.06026"
$ws.Range("D7").Value  = "'0.8082"
$ws.Range("D8").Value  = "'0.9286"
$ws.Range("D9").Value  = "'0.1426"
$ws.Range("D10").Value = "'0.07447"
$ws.Range("D11").Value = "'0.03359"
$ws.Range("D12").Value = "'0.03049"
$ws.Range("D13").Value = "'0.09365"
$ws.Range("D14").Value = "'3.935"
$ws.Range("D15").Value = "'0.001595"
$ws.Range("D16").Value = "'0.04836"

# Row 17 label change (Volume(1h) text)
$ws.Range("E17").Value = "16OneONEWorstin24h"

$ws.Range("D18").Value = "'0.005351"
$ws.Range("D19").Value = "'0.004149"
$ws.Range("D20").Value = "'0.0009883"
$ws.Range("D23").Value = "'6.443"
$ws.Range("D24").Value = "'2.187"
$ws.Range("D40").Value = "'0.03969"

# --- Rows 41-43: coin ranking rotated by one (KickToken/BKEXToken/CEJI) ---
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").Value = "'0.006415"
$ws.Range("E41").Value = "40KickTokenKICK"

$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Value = "'0.1072"
$ws.Range("E42").Value = "41BKEXTokenBKK"

$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D43").Value = "'0.002902"
$ws.Range("E43").Value = "42CEJICEJI"

$ws.Range("D44").Value = "'0.005967"
$ws.Range("D45").Value = "'0.00005207"
$ws.Range("D46").Value = "'0.00000000751"

$ws.Range("D49").Value = "'0.002279"
$ws.Range("E49").Value = "48BOLOBOLO"
$ws.Range("D50").Value = "'0.00002102"
